# Update the "想去人数" (F column) counts on the data sheets.
# Both the "展览" sheet and the "全部类型" sheet carry the same rows of
# event data (rows 2-21), so the same set of F-column updates applies to
# each of them.

$wb = $excel.ActiveWorkbook

$updates = @{
    2  = 267
    3  = 283
    4  = 291
    5  = 844
    6  = 13
    7  = 300
    8  = 8142
    9  = 74
    12 = 109
    15 = 21
    18 = 251
    19 = 703
    20 = 26
    21 = 79
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Cells.Item($row, 6).Value = $updates[$row]
    }
}
